$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update account number (C2) - keep as text value (leading apostrophe
# forces text storage, matching the cell's existing quote-prefixed style)
$ws.Range("C2").Value = "'1010826108"

# Update the transaction date/time (H2)
$ws.Range("H2").Value = "30 jun. 2023, 14:45:36"

# Move the active selection to E7 (matches new sheetView selection)
$ws.Range("E7").Select()
